# DP - create_forecast_ad_hoc - setup main.ipynb
#
# The "מיקום תוכנת תחזית בסיס" / forecast-basic-location row (row 2) is no
# longer needed as an input for the ad-hoc forecast notebook, so it is
# removed entirely - deleting the whole row shifts everything below it up
# by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()

# Leave the selection where the deleted row used to be, matching where the
# cursor lands after removing row 2.
$ws.Range("A2").Select() | Out-Null
